# Etapa7.pptx - "Atualização da aula 5" commit
#
# 1) Slide 3 and Slide 4: class "No" field declarations made private
#    (the "public" modifier on the `conteudo` and `proximoNo` fields
#    becomes "private"; the class declaration itself stays "public").
# 2) Slide 5: reposition the "Imagem 6" picture (move it up).

$p = $ppt.ActivePresentation

function Set-FieldToPrivate {
    # positional params - this runtime's PowerShell interpreter does not
    # reliably bind named (-Param value) arguments to custom functions
    param($Slide, $ShapeIndex, $ParagraphIndex)

    $tr = $Slide.Shapes.Item($ShapeIndex).TextFrame.TextRange
    $para = $tr.Paragraphs($ParagraphIndex)
    $localIdx = $para.Text.IndexOf("public")
    if ($localIdx -ge 0) {
        $globalStart = $para.Start + $localIdx
        $word = $tr.Characters($globalStart, 6)
        $word.Text = "private"
    }
}

# --- Slide 3: class No { public String conteudo; public No proximoNo; ... }
$slide3 = $p.Slides.Item(3)
Set-FieldToPrivate $slide3 4 3   # "    public String conteudo;"
Set-FieldToPrivate $slide3 4 4   # "    public No proximoNo = null;"

# --- Slide 4: class No<T> { public T conteudo; public No proximoNo; ... }
$slide4 = $p.Slides.Item(4)
Set-FieldToPrivate $slide4 4 3   # "    public T conteudo;"
Set-FieldToPrivate $slide4 4 4   # "    public No proximoNo = null;"

# --- Slide 5: move "Imagem 6" picture up (y: 2461450 -> 79756 EMU)
$slide5 = $p.Slides.Item(5)
$img = $slide5.Shapes.Item(12)
$img.Top = 79756 / 914400 * 72
